$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and volume change (E) values
$ws.Range("D2").Value = "26.975.57"
$ws.Range("E2").Value = "  -3.21%  "
$ws.Range("D3").Value = "1.796.46"
$ws.Range("E3").Value = "  -3.33%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -2.75%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4192"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -3.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3567"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -3.82%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07085"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8440"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -4.23%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.11"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.16%  "
$ws.Range("D12").Value = "1.795.70"
$ws.Range("E12").Value = "  -5.32%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.282"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -3.50%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.341"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06758"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.08%  "
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "79.52"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008647"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.94%  "
$ws.Range("E19").Value = "  +0.20%  "
$ws.Range("E20").Value = "  -4.04%  "
$ws.Range("D21").Value = "27.015.73"
$ws.Range("E21").Value = "  -3.81%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.051"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.92"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("D24").Value = "2.015.33"
$ws.Range("E24").Value = "  -5.27%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.934"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.59"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.09"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "4.996"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.91"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.638"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -12.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08957"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7166"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -9.36%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.860"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.289"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -7.42%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.001"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.074"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -9.09%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.073"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.55%  "
$ws.Range("E38").Value = "  -3.59%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05104"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -6.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.1624"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4940"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.570"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -9.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.964"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -12.22%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.018"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -8.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "104.35"
$ws.Range("D45").ClearFormats()
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.18"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.001"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06299"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.4515"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -6.03%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.593"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -4.97%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "61.96"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -5.30%  "
